# Add "Debts" and "Fixed Assets" worksheets (after Kim/Sam), populate their
# header rows, apply bold headers + currency-bold for the money columns, and
# leave "Kim" as the selected/active sheet (matching the target workbook).

$wb = $excel.ActiveWorkbook

# --- create the two new sheets at the end of the tab strip ---------------
$afterSam   = $wb.Worksheets.Item($wb.Worksheets.Count)
$debts      = $wb.Worksheets.Add($null, $afterSam)
$debts.Name = "Debts"

$afterDebts      = $wb.Worksheets.Item($wb.Worksheets.Count)
$fixedAssets     = $wb.Worksheets.Add($null, $afterDebts)
$fixedAssets.Name = "Fixed Assets"

# --- Fixed Assets header row (populated first so its new shared strings
#     "name","type","basis","value","rate","yod","commission" land before
#     the Debts-only strings "term"/"amount") -------------------------------
$fixedAssets.Range("A1").Value = "name"
$fixedAssets.Range("B1").Value = "type"
$fixedAssets.Range("C1").Value = "basis"
$fixedAssets.Range("D1").Value = "value"
$fixedAssets.Range("E1").Value = "rate"
$fixedAssets.Range("F1").Value = "yod"
$fixedAssets.Range("G1").Value = "commission"
$fixedAssets.Range("A1:G1").Font.Bold = $true
$fixedAssets.Range("C1:D1").NumberFormat = """$""#,##0"
$fixedAssets.Rows.Item(1).Select()

# --- Debts header row ------------------------------------------------------
$debts.Range("A1").Value = "name"
$debts.Range("B1").Value = "type"
$debts.Range("C1").Value = "year"
$debts.Range("D1").Value = "term"
$debts.Range("E1").Value = "amount"
$debts.Range("F1").Value = "rate"
$debts.Range("A1:F1").Font.Bold = $true
$debts.Range("E1").NumberFormat = """$""#,##0"
$debts.Rows.Item(1).Select()

# --- restore "Kim" as the active/selected sheet ----------------------------
$kim = $wb.Worksheets.Item("Kim")
$kim.Select()
